$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Update")
Write-Host $ws.Name
